$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" header on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet at the end ---
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# Reuse the bold/centered header style from the Weekly Quantity header row
$wsWeekly.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Reuse the datetime-number-format style from the Weekly Quantity date column
$wsWeekly.Range("A2").Copy()
$ws3.Range("A2:A20").PasteSpecial(-4122)

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

$ws3.Cells.Item(2,1).Value = 45494.99999999999
$ws3.Cells.Item(2,2).Value = 19
$ws3.Cells.Item(2,3).Value = -67.66950355233121
$ws3.Cells.Item(2,4).Value = 100.3159955478662
$ws3.Cells.Item(3,1).Value = 45501.99999999999
$ws3.Cells.Item(3,2).Value = 26
$ws3.Cells.Item(3,3).Value = -49.87345786671792
$ws3.Cells.Item(3,4).Value = 100.6409310031277
$ws3.Cells.Item(4,1).Value = 45515.99999999999
$ws3.Cells.Item(4,2).Value = 38
$ws3.Cells.Item(4,3).Value = -40.3884513316171
$ws3.Cells.Item(4,4).Value = 125.1543256468661
$ws3.Cells.Item(5,1).Value = 45522.99999999999
$ws3.Cells.Item(5,2).Value = 45
$ws3.Cells.Item(5,3).Value = -35.47951298690523
$ws3.Cells.Item(5,4).Value = 126.6427010868983
$ws3.Cells.Item(6,1).Value = 45529.99999999999
$ws3.Cells.Item(6,2).Value = 51
$ws3.Cells.Item(6,3).Value = -25.88191703029135
$ws3.Cells.Item(6,4).Value = 134.1455836590783
$ws3.Cells.Item(7,1).Value = 45536.99999999999
$ws3.Cells.Item(7,2).Value = 58
$ws3.Cells.Item(7,3).Value = -27.67586534479397
$ws3.Cells.Item(7,4).Value = 134.3802074438634
$ws3.Cells.Item(8,1).Value = 45543.99999999999
$ws3.Cells.Item(8,2).Value = 64
$ws3.Cells.Item(8,3).Value = -16.11733200786935
$ws3.Cells.Item(8,4).Value = 150.8042835833479
$ws3.Cells.Item(9,1).Value = 45550.99999999999
$ws3.Cells.Item(9,2).Value = 71
$ws3.Cells.Item(9,3).Value = -8.614536674724873
$ws3.Cells.Item(9,4).Value = 152.5320757798558
$ws3.Cells.Item(10,1).Value = 45557.99999999999
$ws3.Cells.Item(10,2).Value = 77
$ws3.Cells.Item(10,3).Value = -3.766198220653991
$ws3.Cells.Item(10,4).Value = 165.3318796396945
$ws3.Cells.Item(11,1).Value = 45592.99999999999
$ws3.Cells.Item(11,2).Value = 110
$ws3.Cells.Item(11,3).Value = 23.4123174498096
$ws3.Cells.Item(11,4).Value = 193.1771067448674
$ws3.Cells.Item(12,1).Value = 45613.99999999999
$ws3.Cells.Item(12,2).Value = 129
$ws3.Cells.Item(12,3).Value = 43.79705095228192
$ws3.Cells.Item(12,4).Value = 215.6875837923046
$ws3.Cells.Item(13,1).Value = 45620.99999999999
$ws3.Cells.Item(13,2).Value = 135
$ws3.Cells.Item(13,3).Value = 54.56226445931316
$ws3.Cells.Item(13,4).Value = 216.2412718404678
$ws3.Cells.Item(14,1).Value = 45627.99999999999
$ws3.Cells.Item(14,2).Value = 142
$ws3.Cells.Item(14,3).Value = 57.64649546219925
$ws3.Cells.Item(14,4).Value = 229.6070562166537
$ws3.Cells.Item(15,1).Value = 45634.99999999999
$ws3.Cells.Item(15,2).Value = 148
$ws3.Cells.Item(15,3).Value = 63.63803641489062
$ws3.Cells.Item(15,4).Value = 229.4685184460312
$ws3.Cells.Item(16,1).Value = 45641.99999999999
$ws3.Cells.Item(16,2).Value = 155
$ws3.Cells.Item(16,3).Value = 75.07151667455119
$ws3.Cells.Item(16,4).Value = 233.5628794327899
$ws3.Cells.Item(17,1).Value = 45648.99999999999
$ws3.Cells.Item(17,2).Value = 161
$ws3.Cells.Item(17,3).Value = 75.9633178919581
$ws3.Cells.Item(17,4).Value = 247.1942482798513
$ws3.Cells.Item(18,1).Value = 45655.99999999999
$ws3.Cells.Item(18,2).Value = 168
$ws3.Cells.Item(18,3).Value = 87.02751843487454
$ws3.Cells.Item(18,4).Value = 252.1225955050013
$ws3.Cells.Item(19,1).Value = 45662.99999999999
$ws3.Cells.Item(19,2).Value = 174
$ws3.Cells.Item(19,3).Value = 91.34410982895571
$ws3.Cells.Item(19,4).Value = 258.0767051570981
$ws3.Cells.Item(20,1).Value = 45669.99999999999
$ws3.Cells.Item(20,2).Value = 181
$ws3.Cells.Item(20,3).Value = 98.68462439458507
$ws3.Cells.Item(20,4).Value = 261.2993831037251

Write-Output "PO Forecast sheet populated"
